# Weekly refresh of the "Agrícola del Norte S.A. de Arica - Nectarín" sheet.
# Two brand-new price observations (week of D=45015) are inserted at the
# top of the data block (rows 89-90), pushing every existing data row
# down by two (old row 89 -> new row 91, ..., old row 115 -> new row 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 89; this shifts the
# old rows 89..115 down to 91..117 and copies row 88's formatting
# (notably the date style on column D) onto the new rows.
$ws.Rows("89:90").Insert()

$newRows = @(
    @{
        Row = 89
        A = 1; B = "Agrícola del Norte S.A. de Arica"; C = "Arica y Parinacota"
        D = 45015; E = 15; F = "Fruta"; G = 100103
        H = "Frutos de hueso (carozo)"; I = 100103006; J = "Nectarín"
        K = "Artic Sprite"; L = "Segunda"
        M = 250; N = 19000; O = 20000; P = 19500
        Q = "`$/caja 20 kilos granel"
        R = "Región de O'Higgins"
        S = 975; T = 20
    },
    @{
        Row = 90
        A = 1; B = "Agrícola del Norte S.A. de Arica"; C = "Arica y Parinacota"
        D = 45015; E = 15; F = "Fruta"; G = 100103
        H = "Frutos de hueso (carozo)"; I = 100103006; J = "Nectarín"
        K = "August Red"; L = "Segunda"
        M = 270; N = 19000; O = 20000; P = 19500
        Q = "`$/caja 20 kilos granel"
        R = "Región de O'Higgins"
        S = 975; T = 20
    }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("H$r").Value = $entry.H
    $ws.Range("I$r").Value = $entry.I
    $ws.Range("J$r").Value = $entry.J
    $ws.Range("K$r").Value = $entry.K
    $ws.Range("L$r").Value = $entry.L
    $ws.Range("M$r").Value = $entry.M
    $ws.Range("N$r").Value = $entry.N
    $ws.Range("O$r").Value = $entry.O
    $ws.Range("P$r").Value = $entry.P
    $ws.Range("Q$r").Value = $entry.Q
    $ws.Range("R$r").Value = $entry.R
    $ws.Range("S$r").Value = $entry.S
    $ws.Range("T$r").Value = $entry.T
}

Write-Output "Inserted 2 new rows (89-90); last data row is now 117"
